$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the Min/Max header into columns E/F (row 2)
$ws.Range("E2").Value = "Min"
$ws.Range("F2").Value = "Max"

# Add new SPI Transmit time values (37ms -> 10ms) into E4/F4
$ws.Range("E4").Value = "10ms"
$ws.Range("F4").Value = "10ms"

# Update the active selection to match the new focus cell
$ws.Range("G5").Select()
